function New-TagRun {
    param($d, $pos, $tagText, $strike)
    $donor = $d.Content
    $donor.Find.Execute("<m>")
    $donorCopy = $donor.Duplicate
    $donorLen = $donorCopy.End - $donorCopy.Start
    $ins = $d.Range($pos, $pos)
    $ins.FormattedText = $donorCopy.FormattedText
    $newRange = $d.Range($pos, $pos + $donorLen)
    $newRange.Text = $tagText
    if ($strike) {
        $newRange.Font.StrikeThrough = 0
    }
    return $pos + $tagText.Length
}

function New-PlainRun {
    param($d, $pos, $text)
    $donor = $d.Content
    $donor.Find.Execute("Soulder")
    $donorCopy = $donor.Duplicate
    $donorLen = $donorCopy.End - $donorCopy.Start
    $ins = $d.Range($pos, $pos)
    $ins.FormattedText = $donorCopy.FormattedText
    $newRange = $d.Range($pos, $pos + $donorLen)
    $newRange.Text = $text
    return $pos + $text.Length
}

function Add-PaTag {
    param($d, $pos, $word)
    $pos = New-TagRun $d $pos "<pa>" $false
    $pos = New-PlainRun $d $pos $word
    $pos = New-TagRun $d $pos "</pa>" $false
    return $pos
}

$d = $word.ActiveDocument

# ---------- Location 1: paragraph 33 - "Amendiers abricots" ----------
# -> <pa>Amendiers</pa> <pa>abricots</pa>
$p = $d.Paragraphs(33).Range
$p.Find.Execute("Amendiers abricots")
$pos = $p.Start
$p.Text = ""

$pos = Add-PaTag $d $pos "Amendiers"
$pos = New-PlainRun $d $pos " "
$pos = Add-PaTag $d $pos "abricots"

# ---------- Location 2: paragraph 36 - "e pavis mericotons alberges " ----------
# -> "e " <pa>pavis</pa> " " <pa>mericotons</pa> " " <pa>alberges</pa> " "
$p = $d.Paragraphs(36).Range
$p.Find.Execute("e pavis mericotons alberges ")
$pos = $p.Start
$p.Text = ""

$pos = New-PlainRun $d $pos "e "
$pos = Add-PaTag $d $pos "pavis"
$pos = New-PlainRun $d $pos " "
$pos = Add-PaTag $d $pos "mericotons"
$pos = New-PlainRun $d $pos " "
$pos = Add-PaTag $d $pos "alberges"
$pos = New-PlainRun $d $pos " "

# ---------- Location 3: paragraph 36 - "abricots &" ----------
# -> <pa strike=0>abricots</pa> " &"
$p = $d.Paragraphs(36).Range
$p.Find.Execute("abricots &")
$pos = $p.Start
$p.Text = ""

$pos = New-TagRun $d $pos "<pa>" $true
$pos = New-PlainRun $d $pos "abricots"
$pos = New-TagRun $d $pos "</pa>" $false
$pos = New-PlainRun $d $pos " &"

# ---------- Location 4: paragraph 37 - "mieulx hantes sur lamendier en escusson" ----------
# -> "mieulx hantes sur l" <pa>amendier</pa> " en escusson"
$p = $d.Paragraphs(37).Range
$p.Find.Execute("mieulx hantes sur lamendier en escusson")
$pos = $p.Start
$p.Text = ""

$pos = New-PlainRun $d $pos "mieulx hantes sur l"
$pos = Add-PaTag $d $pos "amendier"
$pos = New-PlainRun $d $pos " en escusson"

Write-Output "done"
